$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.452.40'
$ws.Range("E2").Value = '  +2.62%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.854.30'
$ws.Range("E3").Value = '  +2.46%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.34'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.613'
$ws.Range("E6").Value = '  +2.53%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.18'
$ws.Range("E8").Value = '  +15.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.308'
$ws.Range("E9").Value = '  +5.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0693'
$ws.Range("E11").Value = '  +3.59%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.115.81'
$ws.Range("E12").Value = '  +2.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.63'
$ws.Range("E13").Value = '  +2.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.851.13'
$ws.Range("E14").Value = '  +2.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.79'
$ws.Range("E15").Value = '  +8.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.666'
$ws.Range("E16").Value = '  +5.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '35.397.48'
$ws.Range("E17").Value = '  +2.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.11'
$ws.Range("E18").Value = '  +2.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '246.64'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0796'
$ws.Range("E20").Value = '  +2.53%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.29'
$ws.Range("E21").Value = '  +9.38%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.76'
$ws.Range("E22").Value = '  +15.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.20'
$ws.Range("E24").Value = '  -0.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.33'
$ws.Range("E25").Value = '  +0.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.96'
$ws.Range("E26").Value = '  +1.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.91'
$ws.Range("E27").Value = '  +3.37%  '
$ws.Range("E28").Value = '  +2.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.587.71'
$ws.Range("E29").Value = '  +47.66%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.33'
$ws.Range("E31").Value = '  +8.17%  '
$ws.Range("E32").Value = '  +3.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.06'
$ws.Range("E33").Value = '  +2.95%  '
$ws.Range("E34").Value = '  +3.44%  '
$ws.Range("E35").Value = '  +3.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.677'
$ws.Range("E36").Value = '  +3.43%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '90.59'
$ws.Range("E37").Value = '  +11.96%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.09'
$ws.Range("E38").Value = '  +1.09%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.343.76'
$ws.Range("E39").Value = '  -1.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.03'
$ws.Range("E40").Value = '  +9.47%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0195'
$ws.Range("E41").Value = '  +4.15%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.43'
$ws.Range("E42").Value = '  +2.74%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.94'
$ws.Range("E43").Value = '  +7.64%  '
$ws.Range("E44").Value = '  +4.58%  '
$ws.Range("B45").Value = 'MXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.84'
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.45'
$ws.Range("E46").Value = '  +1.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0519'
$ws.Range("E47").Value = '  +3.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.05'
$ws.Range("E48").Value = '  +3.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.015.32'
$ws.Range("E49").Value = '  +2.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '105.15'
$ws.Range("E50").Value = '  +2.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
